$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3200.1667
$ws.Range("I40").Value = 5400.5
$ws.Range("J40").Value = 2100
$ws.Range("K40").Value = 5400.5
$ws.Range("L40").Value = 2100
$ws.Range("M40").Value = -5225.5
$ws.Range("N40").Value = -2450

$ws.Range("H74").Value = 3261.4
$ws.Range("I74").Value = 3261.4
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3261.4
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = -2325.4
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 3261.4
$ws.Range("I77").Value = 3261.4
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 16307
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -11627
$ws.Range("N77").Value = -11627

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("N81").Value = 0
$ws.Range("L81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("N84").Value = 0
$ws.Range("L84").ClearContents()

$ws.Range("H112").Value = 4139.9165
$ws.Range("I112").Value = 1174.5
$ws.Range("J112").Value = 4733
$ws.Range("K112").Value = 3523.5
$ws.Range("L112").Value = 14199
$ws.Range("M112").Value = -2415.5
$ws.Range("N112").Value = -16415

$ws.Range("H116").Value = 2695.889
$ws.Range("I116").Value = 2160.8333
$ws.Range("J116").Value = 3766
$ws.Range("K116").Value = 2160.8333
$ws.Range("L116").Value = 3766
$ws.Range("M116").Value = 1281.1667
$ws.Range("N116").Value = -10650

$ws.Range("H129").Value = 873.1613
$ws.Range("I129").Value = 747.5
$ws.Range("J129").Value = 891.7778
$ws.Range("K129").Value = 2242.5
$ws.Range("L129").Value = 2675.3334
$ws.Range("M129").Value = 2757.5
$ws.Range("N129").Value = -12675.3334

$ws.Range("H132").Value = 6412331.5
$ws.Range("I132").Value = 7577737.5
$ws.Range("K132").Value = 22733212.5
$ws.Range("M132").Value = -22730682.5

$ws.Range("H138").Value = 1696.97
$ws.Range("I138").Value = 686.3333
$ws.Range("J138").Value = 1875.3176
$ws.Range("K138").Value = 2058.9999
$ws.Range("L138").Value = 5625.9528
$ws.Range("M138").Value = 3081.0001
$ws.Range("N138").Value = -15905.9528

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2924.6072
$ws.Range("I32").Value = 3102.9575
$ws.Range("K32").Value = 3102.9575
$ws.Range("M32").Value = -2815.9575

$ws.Range("H132").Value = 2410.8408
$ws.Range("I132").Value = 2171.25
$ws.Range("J132").Value = 3489
$ws.Range("K132").Value = 6513.75
$ws.Range("L132").Value = 10467
$ws.Range("M132").Value = -3983.75
$ws.Range("N132").Value = -15527

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2465.9363
$ws.Range("I86").Value = 2754.7666
$ws.Range("J86").Value = 1956.2354
$ws.Range("K86").Value = 2754.7666
$ws.Range("L86").Value = 1956.2354
$ws.Range("M86").Value = -1631.7666
$ws.Range("N86").Value = -4202.2354

$ws.Range("H89").Value = 2465.9363
$ws.Range("I89").Value = 2754.7666
$ws.Range("J89").Value = 1956.2354
$ws.Range("K89").Value = 13773.833
$ws.Range("L89").Value = 9781.177
$ws.Range("M89").Value = -8157.832999999999
$ws.Range("N89").Value = -21013.177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 5000501
$ws.Range("I12").Value = 668.3333
$ws.Range("K12").Value = 668.3333
$ws.Range("M12").Value = -498.3333

$ws.Range("H19").Value = 325
$ws.Range("I19").Value = 325
$ws.Range("K19").Value = 325
$ws.Range("M19").Value = -155

$ws.Range("H24").Value = 325
$ws.Range("I24").Value = 325
$ws.Range("K24").Value = 325
$ws.Range("M24").Value = -155

$ws.Range("H31").Value = 1064.4906
$ws.Range("I31").Value = 696.6579
$ws.Range("K31").Value = 696.6579
$ws.Range("M31").Value = -401.6579

$ws.Range("H34").Value = 1064.4906
$ws.Range("I34").Value = 696.6579
$ws.Range("K34").Value = 696.6579
$ws.Range("M34").Value = -494.6579

$ws.Range("H58").Value = 902.625
$ws.Range("I58").Value = 696.13336
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 696.13336
$ws.Range("L58").Value = 4000
$ws.Range("M58").Value = -493.13336
$ws.Range("N58").Value = -4406

$ws.Range("H107").Value = 663.0714
$ws.Range("I107").Value = 320
$ws.Range("K107").Value = 320
$ws.Range("M107").Value = 1600

$ws.Range("H132").Value = 6864
$ws.Range("I132").Value = 8101.4116
$ws.Range("J132").Value = 3858.8572
$ws.Range("K132").Value = 24304.2348
$ws.Range("L132").Value = 11576.5716
$ws.Range("M132").Value = -21774.2348
$ws.Range("N132").Value = -16636.5716

$ws.Range("H134").Value = 10102153
$ws.Range("I134").Value = 11905733
$ws.Range("K134").Value = 35717199
$ws.Range("M134").Value = -35714664

$ws.Range("H136").Value = 902.625
$ws.Range("I136").Value = 696.13336
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 2088.40008
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = 461.5999199999997
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1090.909
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2831

$ws.Range("H34").Value = 4168076.2
$ws.Range("J34").Value = 5264772.5
$ws.Range("L34").Value = 15794317.5
$ws.Range("N34").Value = -15794485.5

$ws.Range("H87").Value = 1626.4
$ws.Range("J87").Value = 3316
$ws.Range("L87").Value = 9948
$ws.Range("N87").Value = -12444

$ws.Range("H90").Value = 1626.4
$ws.Range("J90").Value = 3316
$ws.Range("L90").Value = 29844
$ws.Range("N90").Value = -42324

$ws.Range("H107").Value = 13617.625
$ws.Range("J107").Value = 26100.75
$ws.Range("L107").Value = 78302.25
$ws.Range("N107").Value = -82142.25

$ws.Range("H131").Value = 12501176
$ws.Range("I131").Value = 111111790
$ws.Range("J131").Value = 1238.5916
$ws.Range("K131").Value = 333335370
$ws.Range("L131").Value = 3715.7748
$ws.Range("M131").Value = -333330330
$ws.Range("N131").Value = -13795.7748

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("N125").Value = 0
$ws.Range("L125").ClearContents()

$ws.Range("H132").Value = 3984.9092
$ws.Range("I132").Value = 4220.5
$ws.Range("J132").Value = 3702.2
$ws.Range("K132").Value = 12661.5
$ws.Range("L132").Value = 11106.6
$ws.Range("M132").Value = -10131.5
$ws.Range("N132").Value = -16166.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 3127
$ws.Range("I31").Value = 1907.5
$ws.Range("J31").Value = 3736.75
$ws.Range("K31").Value = 1907.5
$ws.Range("L31").Value = 3736.75
$ws.Range("M31").Value = -1659.5
$ws.Range("N31").Value = -4232.75

$ws.Range("H46").Value = 6160
$ws.Range("I46").Value = 542
$ws.Range("J46").Value = 8969
$ws.Range("K46").Value = 542
$ws.Range("L46").Value = 8969
$ws.Range("M46").Value = -354
$ws.Range("N46").Value = -9345

$ws.Range("H55").Value = 790.8182
$ws.Range("I55").Value = 233.66667
$ws.Range("J55").Value = 999.75
$ws.Range("K55").Value = 233.66667
$ws.Range("L55").Value = 999.75
$ws.Range("M55").Value = -60.66667000000001
$ws.Range("N55").Value = -1345.75

$ws.Range("H136").Value = 6926.1763
$ws.Range("I136").Value = 7569.6665
$ws.Range("J136").Value = 2100
$ws.Range("K136").Value = 22708.9995
$ws.Range("L136").Value = 6300
$ws.Range("M136").Value = -20158.9995
$ws.Range("N136").Value = -11400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 558.7646999999999
$ws.Range("I107").Value = 485.42856
$ws.Range("J107").Value = 901
$ws.Range("K107").Value = 1456.28568
$ws.Range("L107").Value = 2703
$ws.Range("M107").Value = 463.71432
$ws.Range("N107").Value = -6543

$ws.Range("H131").Value = 74906
$ws.Range("J131").Value = 74906
$ws.Range("L131").Value = 74906
$ws.Range("N131").Value = -84986

$ws.Range("H136").Value = 569.1539
$ws.Range("I136").Value = 281.125
$ws.Range("K136").Value = 843.375
$ws.Range("M136").Value = 1706.625
